$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New normalized values (row -> column letter -> value)
$values = @{
    2 = @{ B = -144.9629974365234;  C = -0.1087;                  D = -0.1348000019788742;  E = 0.5033000111579895;  F = -0.1348000019788742 }
    3 = @{ B = -155.5372924804688;  C = -0.1244;                  D = -0.1348;               E = 0.3291000127792358;  F = -0.1348000019788742 }
    4 = @{ B = -114.5098037719727;  C = -0.1338;                  D = -0.1348;               E = 0.2471999973058701;  F = -0.1348000019788742 }
    5 = @{ B = 130.3462066650391;   C = 0.1375;                   D = 0.1698;                E = 0.6880000233650208;  F = -0.1348000019788742 }
    6 = @{ B = 228.21240234375;     C = 0.1805;                   D = 0.1162;                E = 0.8651999831199646;  F = -0.1348000019788742 }
    7 = @{ B = 173.0357971191406;   C = 0.1743;                   D = 0.1578000038862228;   E = 0.5821999907493591;  F = -0.1348000019788742 }
    8 = @{ B = 258.4993896484375;   C = 0.2316;                   D = 0.2115;                E = 0.5821999907493591;  F = -0.1348000019788742 }
    9 = @{ B = 375.0837097167969;   C = 0.048;                    D = 0;                     E = 0.8651999831199646;  F = -0.1348000019788742 }
}

foreach ($row in $values.Keys) {
    foreach ($col in $values[$row].Keys) {
        $ws.Range("$col$row").Value = $values[$row][$col]
    }
}
